$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift the existing rating-date columns (B:E = Jun_17/Jun_15/Jun_13/Jun_10)
# three places to the right (-> E:H) to make room for the two new "Jun_26"
# columns and the new "Jun_27" column. Using a manual copy (rather than
# Range.Insert) keeps the pre-existing column width metadata on columns
# C/D/E untouched, matching the target layout.
for ($r = 1; $r -le 27; $r++) {
    $valB = $ws.Cells.Item($r, 2).Value2
    $valC = $ws.Cells.Item($r, 3).Value2
    $valD = $ws.Cells.Item($r, 4).Value2
    $valE = $ws.Cells.Item($r, 5).Value2
    $ws.Cells.Item($r, 8).Value = $valE
    $ws.Cells.Item($r, 7).Value = $valD
    $ws.Cells.Item($r, 6).Value = $valC
    $ws.Cells.Item($r, 5).Value = $valB
}

# New "Jun_26" rating column pair (two analysts rated on Jun_26)
$ws.Range("C1").Value = "Jun_26"
$ws.Range("D1").Value = "Jun_26"
for ($r = 2; $r -le 27; $r++) {
    $ws.Cells.Item($r, 3).Value = "UN"
    $ws.Cells.Item($r, 4).Value = "UN"
}

# New analyst rows added to the watch list
$ws.Range("A28").Value = "Benchmark"
$ws.Range("B28").Value = "UN"
$ws.Range("C28").Value = "UN"
$ws.Range("D28").Value = "UN"

$ws.Range("A29").Value = "Evercore ISI"
$ws.Range("B29").Value = "UN"
$ws.Range("C29").Value = "UN"
$ws.Range("D29").Value = "UN"

# Newest "Jun_27" rating column (inserted at the front, column B)
$ws.Range("B1").Value = "Jun_27"
for ($r = 2; $r -le 27; $r++) {
    $ws.Cells.Item($r, 2).Value = "UN"
}

# Give the newly-used columns (F,G,H) the same display width as the
# existing C:E columns (8 characters).
$ws.Columns(6).ColumnWidth = 7.1
$ws.Columns(7).ColumnWidth = 7.1
$ws.Columns(8).ColumnWidth = 7.1
